$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing quantity values (column F)
$ws.Range("F2").Value = -233
$ws.Range("F3").Value = -441
$ws.Range("F4").Value = -706

# Row 5: location changes from DC_002 to DC_001, quantity and horizon_days change
$ws.Range("B5").Value = "DC_001"
$ws.Range("F5").Value = -103
$ws.Range("H5").Value = 4

# Add new row 6: MAT_B / DC_002
$ws.Range("A6").Value = "MAT_B"
$ws.Range("B6").Value = "DC_002"
$ws.Range("C6").Value = 45295
$ws.Range("D6").Value = "Distribution Demand - Forecast"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = -33
$ws.Range("G6").Value = 45294
$ws.Range("H6").Value = 1

# Add new row 7: MAT_B / PLANT_001
$ws.Range("A7").Value = "MAT_B"
$ws.Range("B7").Value = "PLANT_001"
$ws.Range("C7").Value = 45295
$ws.Range("D7").Value = "Distribution Demand - Forecast"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -103
$ws.Range("G7").Value = 45294
$ws.Range("H7").Value = 1

# Apply the same date number format as the other date cells (numFmtId 165)
$ws.Range("C6:C7").NumberFormat = $ws.Range("C5").NumberFormat
$ws.Range("G6:G7").NumberFormat = $ws.Range("G5").NumberFormat
